$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($null -ne $val -and $val.Length -gt 2 -and $val.EndsWith("16")) {
        $cell.Value2 = $val.Substring(0, $val.Length - 2)
    }
}
